# Applies the "B5 netwerken / overzicht" edit:
#   P2: netwerkonderdelen          -> P2: topologie / " en onderdelen"  (2 runs)
#   P3: IP (v4 en v6)              -> removed entirely
#   P / 4 / : cloud                -> P / 3 / : cloud   (renumbered)
#   P5 / bookmark / : protocollen  -> P4 / bookmark / : protocollen (renumbered)

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, $innerXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($pkg)
}

# Locate the paragraphs by their (unique) text rather than trusting fixed
# indices, so the script is robust to how the document is currently laid out.
function Find-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs($i)
        # Paragraph.Range.Text includes the trailing paragraph mark (CR) -
        # trim it before comparing against the plain needle text.
        if ($candidate.Range.Text.TrimEnd() -eq $needle) {
            return $candidate
        }
    }
    return $null
}

# --- 1) "P2: netwerkonderdelen" -> "P2: topologie" + " en onderdelen" (two runs) ---
$p2 = Find-ParagraphByText "P2: netwerkonderdelen"
if ($p2 -eq $null) { throw "Could not locate paragraph 'P2: netwerkonderdelen'" }
Set-ParagraphXml $p2 '<w:p w:rsidR="001C081A" w:rsidRDefault="001C081A"><w:r><w:t>P2: topologie</w:t></w:r><w:r><w:t xml:space="preserve"> en onderdelen</w:t></w:r></w:p>'

# --- 2) Remove the "P3: IP (v4 en v6)" paragraph entirely ---
$pIp = Find-ParagraphByText "P3: IP (v4 en v6)"
if ($pIp -eq $null) { throw "Could not locate paragraph 'P3: IP (v4 en v6)'" }
$delRange = $d.Range($pIp.Range.Start, $pIp.Range.End)
$delRange.Delete()

# --- 3) "P" / "4" / ": cloud" -> "P" / "3" / ": cloud" ---
$pCloud = Find-ParagraphByText "P4: cloud"
if ($pCloud -eq $null) { throw "Could not locate paragraph 'P4: cloud'" }
Set-ParagraphXml $pCloud '<w:p w:rsidR="00174042" w:rsidRDefault="00174042"><w:r><w:t>P</w:t></w:r><w:r><w:t>3</w:t></w:r><w:r><w:t>: cloud</w:t></w:r></w:p>'

# --- 4) "P5" -> "P4" (keep bookmark + ": protocollen" run untouched) ---
$d.Content.Find.Execute("P5", $true, $false, $false, $false, $false, $true, 1, $false, "P4", 2)
